$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @{ D = <new price text>; E = <new volume text> }
# Only columns that actually changed for a given coin row are listed.
$updates = @{
    2 = @{ D='30.279.03'; E='  +1.99%  ' }
    3 = @{ D='1.895.86'; E='  -0.72%  ' }
    4 = @{ E='  -0.15%  ' }
    5 = @{ D='324.81'; E='  +2.77%  ' }
    6 = @{ E='  -0.13%  ' }
    7 = @{ D='0.5175'; E='  +0.12%  ' }
    8 = @{ D='0.4016'; E='  +1.08%  ' }
    9 = @{ D='0.08411'; E='  -1.01%  ' }
    10 = @{ D='42.74'; E='  +0.14%  ' }
    11 = @{ D='1.115' }
    12 = @{ D='23.23'; E='  +10.99%  ' }
    13 = @{ D='6.435'; E='  +2.02%  ' }
    14 = @{ D='1.909.95'; E='  +0.17%  ' }
    15 = @{ D='7.328'; E='  -0.32%  ' }
    16 = @{ D='1.001'; E='  -0.15%  ' }
    17 = @{ D='94.19'; E='  +0.85%  ' }
    18 = @{ E='  -0.72%  ' }
    19 = @{ D='0.06643'; E='  -1.69%  ' }
    20 = @{ E='  +1.43%  ' }
    21 = @{ E='  -0.09%  ' }
    22 = @{ D='5.952'; E='  -1.43%  ' }
    23 = @{ D='30.260.11'; E='  +1.88%  ' }
    24 = @{ D='11.30'; E='  +0.60%  ' }
    25 = @{ E='  +0.81%  ' }
    26 = @{ D='2.110.72'; E='  -0.63%  ' }
    27 = @{ D='21.64'; E='  +3.16%  ' }
    28 = @{ D='161.80'; E='  +1.65%  ' }
    29 = @{ D='2.355'; E='  -3.83%  ' }
    30 = @{ D='129.19'; E='  +0.78%  ' }
    31 = @{ D='1.090'; E='  +1.11%  ' }
    32 = @{ E='  +0.09%  ' }
    33 = @{ D='6.095'; E='  -1.53%  ' }
    34 = @{ D='3.748'; E='  +2.24%  ' }
    35 = @{ E='  -0.11%  ' }
    36 = @{ D='0.06558'; E='  -1.07%  ' }
    37 = @{ D='5.291'; E='  +1.42%  ' }
    38 = @{ E='  -0.33%  ' }
    39 = @{ D='1.222'; E='  -1.27%  ' }
    40 = @{ D='11.80'; E='  +4.16%  ' }
    41 = @{ D='8.751'; E='  -3.75%  ' }
    42 = @{ D='0.6497'; E='  -0.87%  ' }
    43 = @{ D='1.232'; E='  -0.64%  ' }
    44 = @{ D='0.6098'; E='  -0.27%  ' }
    45 = @{ D='13.27'; E='  +0.70%  ' }
    46 = @{ D='3.689'; E='  +0.24%  ' }
    47 = @{ D='2.056'; E='  -0.42%  ' }
    48 = @{ D='1.236'; E='  -0.12%  ' }
    49 = @{ D='124.66'; E='  +0.34%  ' }
    50 = @{ D='1.162'; E='  +0.22%  ' }
    51 = @{ D='79.14'; E='  +1.12%  ' }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"   # keep "1.30"-style text from collapsing to a number
        $cell.Value = $rowData["D"]
    }
    if ($rowData.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $rowData["E"]
    }
}
